$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 'Vega Monumental Concepción'
$ws.Range("C2").Value = 'Bíobío'
$ws.Range("D2").Value = 44477
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 300000000
$ws.Range("G2").Value = 'Espárragos'
$ws.Range("H2").Value = 'Sin especificar'
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = 1460
$ws.Range("N2").Value = '$/kilo'
$ws.Range("O2").Value = 'Provincia de Linares'
$ws.Range("P2").Value = 1460
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 'Hortaliza'

# Row 3
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 'Vega Monumental Concepción'
$ws.Range("C3").Value = 'Bíobío'
$ws.Range("D3").Value = 44860
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 300000000
$ws.Range("G3").Value = 'Espárragos'
$ws.Range("H3").Value = 'Sin especificar'
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 1100
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = 1609
$ws.Range("N3").Value = '$/kilo'
$ws.Range("O3").Value = 'Provincia de Linares'
$ws.Range("P3").Value = 1609
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 'Hortaliza'

# Row 4
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = 'Vega Monumental Concepción'
$ws.Range("C4").Value = 'Bíobío'
$ws.Range("D4").Value = 44496
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 300000000
$ws.Range("G4").Value = 'Espárragos'
$ws.Range("H4").Value = 'Sin especificar'
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 1773
$ws.Range("N4").Value = '$/paquete'
$ws.Range("O4").Value = 'Provincia de Linares'
$ws.Range("P4").Value = 1773
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 'Hortaliza'

# Row 5
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = 'Vega Monumental Concepción'
$ws.Range("C5").Value = 'Bíobío'
$ws.Range("D5").Value = 44868
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 300000000
$ws.Range("G5").Value = 'Espárragos'
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1300
$ws.Range("M5").Value = 1250
$ws.Range("N5").Value = '$/kilo'
$ws.Range("O5").Value = 'Región del Maule'
$ws.Range("P5").Value = 1250
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 'Hortaliza'

# Row 6
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = 'Vega Monumental Concepción'
$ws.Range("C6").Value = 'Bíobío'
$ws.Range("D6").Value = 44868
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 300000000
$ws.Range("G6").Value = 'Espárragos'
$ws.Range("H6").Value = 'Sin especificar'
$ws.Range("I6").Value = 'Segunda'
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 1000
$ws.Range("N6").Value = '$/kilo'
$ws.Range("O6").Value = 'Región del Maule'
$ws.Range("P6").Value = 1000
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 'Hortaliza'

# Row 7
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = 'Vega Monumental Concepción'
$ws.Range("C7").Value = 'Bíobío'
$ws.Range("D7").Value = 44881
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 300000000
$ws.Range("G7").Value = 'Espárragos'
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 2600
$ws.Range("L7").Value = 2700
$ws.Range("M7").Value = 2650
$ws.Range("N7").Value = '$/kilo'
$ws.Range("O7").Value = 'Provincia de Linares'
$ws.Range("P7").Value = 2650
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 'Hortaliza'

# Row 8
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = 'Vega Monumental Concepción'
$ws.Range("C8").Value = 'Bíobío'
$ws.Range("D8").Value = 44881
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 300000000
$ws.Range("G8").Value = 'Espárragos'
$ws.Range("H8").Value = 'Sin especificar'
$ws.Range("I8").Value = 'Segunda'
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 2400
$ws.Range("L8").Value = 2400
$ws.Range("M8").Value = 2400
$ws.Range("N8").Value = '$/kilo'
$ws.Range("O8").Value = 'Provincia de Linares'
$ws.Range("P8").Value = 2400
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 'Hortaliza'

# Row 9
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = 'Vega Monumental Concepción'
$ws.Range("C9").Value = 'Bíobío'
$ws.Range("D9").Value = 45245
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 300000000
$ws.Range("G9").Value = 'Espárragos'
$ws.Range("H9").Value = 'Sin especificar'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 1800
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = 1900
$ws.Range("N9").Value = '$/kilo'
$ws.Range("O9").Value = 'Provincia de Linares'
$ws.Range("P9").Value = 1900
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 'Hortaliza'

# Row 10
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = 'Vega Monumental Concepción'
$ws.Range("C10").Value = 'Bíobío'
$ws.Range("D10").Value = 44875
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 300000000
$ws.Range("G10").Value = 'Espárragos'
$ws.Range("H10").Value = 'Sin especificar'
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 1600
$ws.Range("M10").Value = 1550
$ws.Range("N10").Value = '$/kilo'
$ws.Range("O10").Value = 'Provincia de Linares'
$ws.Range("P10").Value = 1550
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 'Hortaliza'

# Row 11
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = 'Vega Monumental Concepción'
$ws.Range("C11").Value = 'Bíobío'
$ws.Range("D11").Value = 44519
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 300000000
$ws.Range("G11").Value = 'Espárragos'
$ws.Range("H11").Value = 'Sin especificar'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 1200
$ws.Range("L11").Value = 1300
$ws.Range("M11").Value = 1240
$ws.Range("N11").Value = '$/kilo'
$ws.Range("O11").Value = 'Provincia de Linares'
$ws.Range("P11").Value = 1240
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 'Hortaliza'

# Row 12
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 'Vega Monumental Concepción'
$ws.Range("C12").Value = 'Bíobío'
$ws.Range("D12").Value = 44489
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = 300000000
$ws.Range("G12").Value = 'Espárragos'
$ws.Range("H12").Value = 'Sin especificar'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 1400
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1450
$ws.Range("N12").Value = '$/kilo'
$ws.Range("O12").Value = 'Provincia de Linares'
$ws.Range("P12").Value = 1450
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 'Hortaliza'

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 'Vega Monumental Concepción'
$ws.Range("C13").Value = 'Bíobío'
$ws.Range("D13").Value = 45203
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 300000000
$ws.Range("G13").Value = 'Espárragos'
$ws.Range("H13").Value = 'Sin especificar'
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 1400
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1450
$ws.Range("N13").Value = '$/kilo'
$ws.Range("O13").Value = 'Provincia de Linares'
$ws.Range("P13").Value = 1450
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 'Hortaliza'

# Row 14
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = 'Vega Monumental Concepción'
$ws.Range("C14").Value = 'Bíobío'
$ws.Range("D14").Value = 44468
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 300000000
$ws.Range("G14").Value = 'Espárragos'
$ws.Range("H14").Value = 'Verde'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 1800
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 1920
$ws.Range("N14").Value = '$/kilo'
$ws.Range("O14").Value = 'Provincia de Linares'
$ws.Range("P14").Value = 1920
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 'Hortaliza'

# Row 15
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = 'Vega Monumental Concepción'
$ws.Range("C15").Value = 'Bíobío'
$ws.Range("D15").Value = 45246
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 300000000
$ws.Range("G15").Value = 'Espárragos'
$ws.Range("H15").Value = 'Sin especificar'
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 2000
$ws.Range("N15").Value = '$/kilo'
$ws.Range("O15").Value = 'Provincia de Linares'
$ws.Range("P15").Value = 2000
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 'Hortaliza'

# Row 16
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = 'Vega Monumental Concepción'
$ws.Range("C16").Value = 'Bíobío'
$ws.Range("D16").Value = 45246
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 300000000
$ws.Range("G16").Value = 'Espárragos'
$ws.Range("H16").Value = 'Sin especificar'
$ws.Range("I16").Value = 'Segunda'
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = 1500
$ws.Range("N16").Value = '$/kilo'
$ws.Range("O16").Value = 'Provincia de Linares'
$ws.Range("P16").Value = 1500
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 'Hortaliza'

# Row 17
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = 'Vega Monumental Concepción'
$ws.Range("C17").Value = 'Bíobío'
$ws.Range("D17").Value = 45230
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 300000000
$ws.Range("G17").Value = 'Espárragos'
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = 1500
$ws.Range("N17").Value = '$/kilo'
$ws.Range("O17").Value = 'Provincia de Linares'
$ws.Range("P17").Value = 1500
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = 'Hortaliza'

# Row 18
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = 'Vega Monumental Concepción'
$ws.Range("C18").Value = 'Bíobío'
$ws.Range("D18").Value = 44511
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 300000000
$ws.Range("G18").Value = 'Espárragos'
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 1300
$ws.Range("L18").Value = 1400
$ws.Range("M18").Value = 1350
$ws.Range("N18").Value = '$/kilo'
$ws.Range("O18").Value = 'Provincia de Linares'
$ws.Range("P18").Value = 1350
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 'Hortaliza'

# Row 19
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = 'Vega Monumental Concepción'
$ws.Range("C19").Value = 'Bíobío'
$ws.Range("D19").Value = 44545
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 300000000
$ws.Range("G19").Value = 'Espárragos'
$ws.Range("H19").Value = 'Sin especificar'
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 550
$ws.Range("K19").Value = 1700
$ws.Range("L19").Value = 1800
$ws.Range("M19").Value = 1755
$ws.Range("N19").Value = '$/kilo'
$ws.Range("O19").Value = 'Provincia de Linares'
$ws.Range("P19").Value = 1755
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 'Hortaliza'

# Row 20
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = 'Vega Monumental Concepción'
$ws.Range("C20").Value = 'Bíobío'
$ws.Range("D20").Value = 44839
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 300000000
$ws.Range("G20").Value = 'Espárragos'
$ws.Range("H20").Value = 'Sin especificar'
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 500
$ws.Range("K20").Value = 1700
$ws.Range("L20").Value = 1800
$ws.Range("M20").Value = 1760
$ws.Range("N20").Value = '$/kilo'
$ws.Range("O20").Value = 'Provincia de Linares'
$ws.Range("P20").Value = 1760
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = 'Hortaliza'

# Row 21
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = 'Vega Monumental Concepción'
$ws.Range("C21").Value = 'Bíobío'
$ws.Range("D21").Value = 44526
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 300000000
$ws.Range("G21").Value = 'Espárragos'
$ws.Range("H21").Value = 'Sin especificar'
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 1600
$ws.Range("M21").Value = 1550
$ws.Range("N21").Value = '$/kilo'
$ws.Range("O21").Value = 'Provincia de Linares'
$ws.Range("P21").Value = 1550
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = 'Hortaliza'

# Row 22
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = 'Vega Monumental Concepción'
$ws.Range("C22").Value = 'Bíobío'
$ws.Range("D22").Value = 45202
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = 300000000
$ws.Range("G22").Value = 'Espárragos'
$ws.Range("H22").Value = 'Verde'
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 1600
$ws.Range("L22").Value = 1600
$ws.Range("M22").Value = 1600
$ws.Range("N22").Value = '$/kilo'
$ws.Range("O22").Value = 'Provincia de Linares'
$ws.Range("P22").Value = 1600
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = 'Hortaliza'

# Row 23
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = 'Vega Monumental Concepción'
$ws.Range("C23").Value = 'Bíobío'
$ws.Range("D23").Value = 44510
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 300000000
$ws.Range("G23").Value = 'Espárragos'
$ws.Range("H23").Value = 'Sin especificar'
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 600
$ws.Range("K23").Value = 1300
$ws.Range("L23").Value = 1400
$ws.Range("M23").Value = 1350
$ws.Range("N23").Value = '$/kilo'
$ws.Range("O23").Value = 'Provincia de Linares'
$ws.Range("P23").Value = 1350
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = 'Hortaliza'

# Row 24
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = 'Vega Monumental Concepción'
$ws.Range("C24").Value = 'Bíobío'
$ws.Range("D24").Value = 44876
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 300000000
$ws.Range("G24").Value = 'Espárragos'
$ws.Range("H24").Value = 'Sin especificar'
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 350
$ws.Range("K24").Value = 1500
$ws.Range("L24").Value = 1600
$ws.Range("M24").Value = 1557
$ws.Range("N24").Value = '$/kilo'
$ws.Range("O24").Value = 'Provincia de Linares'
$ws.Range("P24").Value = 1557
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = 'Hortaliza'

# Row 25
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = 'Vega Monumental Concepción'
$ws.Range("C25").Value = 'Bíobío'
$ws.Range("D25").Value = 44524
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 300000000
$ws.Range("G25").Value = 'Espárragos'
$ws.Range("H25").Value = 'Sin especificar'
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 1600
$ws.Range("M25").Value = 1550
$ws.Range("N25").Value = '$/kilo'
$ws.Range("O25").Value = 'Provincia de Talca'
$ws.Range("P25").Value = 1550
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = 'Hortaliza'
